# Generate Report for Handoff
# Replace the first handoff record's GUID-named files (1b1858c3-...) with
# 49dbc38d-2320-4cc0-81ef-5fd882349c85, and the second record's GUID-named
# files (c418e212-...) with ffff2d07819d-df09-45b7-8487-cdf5be329d8e.
# Status moves from "Handed back: in sync with en-US" to "Ready for handoff",
# handoff timestamps advance, target/handback info is cleared (new handoff,
# not yet handed back) and "Content Duplicate" flips to True for the second
# record.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.md"
$ov.Range("A3").Value = "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md"

$ov.Range("B2").Value = "e2e\49dbc38d-2320-4cc0-81ef-5fd882349c85.md"
$ov.Range("B3").Value = "e2e\ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md"

foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') { $hl.TextToDisplay = "e2e\49dbc38d-2320-4cc0-81ef-5fd882349c85.md" }
    if ($addr -eq '$B$3') { $hl.TextToDisplay = "e2e\ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md" }
}

$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-09-02 11:13:55"

$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-02 11:13:55"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.md"
$zh.Range("A3").Value = "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md"

foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "49dbc38d-2320-4cc0-81ef-5fd882349c85.md" }
    if ($addr -eq '$A$3') { $hl.TextToDisplay = "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md" }
}

$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("C3").Value = "Ready for handoff"

$zh.Range("G2").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.193e86ffbf271f42e1c4866a2fdea8c483367546.zh-cn.xlf"
$zh.Range("G3").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.193e86ffbf271f42e1c4866a2fdea8c483367546.zh-cn.xlf"

$zh.Range("H2").Value = "2016-09-02 11:13:51"
$zh.Range("H3").Value = "2016-09-02 11:13:51"

$zh.Range("K2").Value = "0001-01-01 00:00:00"
$zh.Range("K3").Value = "0001-01-01 00:00:00"

# Second record's "Content Duplicate" flips False -> True (keep as text).
$zh.Range("F3").Value = "'True"
$zh.Range("F3").Style = "Normal"

# Drop the now-stale "Latest Target File" / "Latest Handback File" links -
# this handoff hasn't been handed back yet.
foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if (($addr -eq '$I$2') -or ($addr -eq '$I$3')) { $hl.Delete() }
}
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Value = ""
$zh.Range("J3").Value = ""
$zh.Range("I3").Value = ""
$zh.Range("I3").Style = "Normal"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.md"
$de.Range("A3").Value = "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md"

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "49dbc38d-2320-4cc0-81ef-5fd882349c85.md" }
    if ($addr -eq '$A$3') { $hl.TextToDisplay = "ffff2d07819d-df09-45b7-8487-cdf5be329d8e.md" }
}

$de.Range("C2").Value = "Ready for handoff"
$de.Range("C3").Value = "Ready for handoff"

$de.Range("G2").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.193e86ffbf271f42e1c4866a2fdea8c483367546.de-de.xlf"
$de.Range("G3").Value = "49dbc38d-2320-4cc0-81ef-5fd882349c85.193e86ffbf271f42e1c4866a2fdea8c483367546.de-de.xlf"

$de.Range("H2").Value = "2016-09-02 11:13:55"
$de.Range("H3").Value = "2016-09-02 11:13:55"

$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Range("K3").Value = "0001-01-01 00:00:00"

# Second record's "Content Duplicate" flips False -> True (keep as text).
$de.Range("F3").Value = "'True"
$de.Range("F3").Style = "Normal"

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if (($addr -eq '$I$2') -or ($addr -eq '$I$3')) { $hl.Delete() }
}
$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("J2").Value = ""
$de.Range("J3").Value = ""
$de.Range("I3").Value = ""
$de.Range("I3").Style = "Normal"
